# Cost Budget Estimate.xlsx update:
# The "Quality Assurance" hourly rate (C7) was revised from 100 to 150.
# This ripples through D7 (=C7*B7), F7 (=D7*E7) and the overall Budget
# total in B1 (=SUM(F6:F9,F12:F20)) automatically via formula recalculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Value = 150

# Leave the cursor where the user would land after editing C7 and
# pressing Enter (one row down).
$ws.Range("C8").Select()
